# Games page in xls updated to reflect .sql
# Rebuild the "Game" worksheet: new headers (GameID, GameName, " ESRB_ID",
# Multiplayer, IsOnline) and refreshed game rows/flags, matching the
# regenerated SQL export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Game")

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "GameID"
$ws.Range("B1").Value = "GameName"
$ws.Range("C1").Value = " ESRB_ID"
$ws.Range("D1").Value = "Multiplayer"
$ws.Range("E1").Value = "IsOnline"

# --- Data rows ---------------------------------------------------------
# Each tuple: GameID, GameName, ESRB_ID, Multiplayer, IsOnline
# NOTE: GameName values genuinely start with a literal apostrophe
# (e.g. 'Final Fantasy I'). In Excel a *leading* apostrophe is the
# "treat as text" entry prefix, so it must be doubled ('') when typed
# in order to land in the cell as a single literal leading quote.
$rows = @(
    @(1, "''Final Fantasy I'", 2, 0, 0),
    @(2, "''Final Fantasy II'", 3, 0, 0),
    @(3, "''Final Fantasy III'", 1, 0, 0),
    @(4, "''Final Fantasy IV'", 2, 0, 0),
    @(5, "''Final Fantasy V'", 3, 0, 0),
    @(6, "''Final Fantasy VI'", 2, 0, 0),
    @(7, "''Final Fantasy VII'", 3, 0, 0),
    @(8, "''Final Fantasy VIII'", 3, 0, 0),
    @(9, "''Final FantasyIX'", 3, 0, 0),
    @(10, "''Final Fantasy X'", 3, 0, 0),
    @(11, "''Final Fantasy XI'", 3, 0, 0),
    @(12, "''Final Fantasy XII'", 3, 0, 0),
    @(13, "''Final Fantasy XIV'", 3, 0, 1),
    @(14, "''Final Fantasy XV'", 3, 0, 0),
    @(15, "''Missile Command'", 1, 0, 0),
    @(16, "''Crash Bandicoot'", 5, 0, 0),
    @(17, "''Star Wars Battlefront'", 3, 1, 0),
    @(18, "''Halo: Combat Evolved'", 4, 1, 1),
    @(19, "''Halo 2'", 4, 1, 1),
    @(20, "''Halo 3'", 4, 1, 1),
    @(21, "''Halo 3: ODST'", 4, 1, 1),
    @(22, "''Halo Reach'", 4, 1, 1),
    @(23, "''Halo 4'", 4, 1, 1),
    @(24, "''Halo 5: Guardians'", 3, 1, 1),
    @(25, "''Halo Wars'", 3, 1, 1),
    @(26, "''Halo Wars 2'", 3, 1, 1),
    @(27, "''Guitar Hero: On Tour'", 2, 1, 1),
    @(28, "''Monster Hunter World'", 3, 1, 1),
    @(29, "''Mario Kart 64'", 1, 1, 0),
    @(30, "''Mario Kart Wii'", 1, 1, 1),
    @(31, "''Left 4 Dead'", 4, 1, 1),
    @(32, "''Elder Scrolls: Arena'", 4, 0, 0),
    @(33, "''Elder Scrolls II: Daggerfall'", 4, 0, 0),
    @(34, "''Elder Scrolls III: Marrowind'", 4, 0, 0),
    @(35, "''Elder Scrolls IV: Oblivion'", 4, 0, 0),
    @(36, "''Elder Scrolls V: Skyrim'", 4, 0, 0),
    @(37, "''Elder Scrolls Online'", 4, 1, 1),
    @(38, "''Destiny'", 3, 1, 1),
    @(39, "''Destiny 2'", 3, 1, 1),
    @(40, "''Super Smash Bros. Brawl'", 3, 1, 1),
    @(41, "''Sonic The Hedgehog'", 1, 1, 0),
    @(42, "''Legend of Zelda'", 1, 0, 0),
    @(43, "''Legend of Zelda: Ocarina of Time'", 1, 0, 0),
    @(44, "''Legend of Zelda: Majoras Mask'", 1, 0, 0),
    @(45, "''Legend of Zelda Links Awakening'", 1, 0, 0),
    @(46, "''Donkey Kong 64'", 1, 1, 0),
    @(47, "''Mario 64'", 1, 0, 0),
    @(48, "''Last of Us'", 4, 1, 1),
    @(49, "''Crazy Taxi'", 3, 0, 0),
    @(50, "''Civilizations'", 2, 1, 1)
)

foreach ($row in $rows) {
    $r = $row[0] + 1
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# --- Selection ----------------------------------------------------------
$ws.Range("A1:E51").Select()
